$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.251.47'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.507.56'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.66%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.543'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.512.14'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.76%  '
$ws.Range('E10').Value = '  +2.13%  '
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.30%  '
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.59'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.97%  '
$ws.Range('E15').Value = '  +0.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.961.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.984.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.521.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.67'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '331.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.24'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.68%  '
$ws.Range('E23').Value = '  +18.84%  '
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.82'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '636.36'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +12.08%  '
$ws.Range('E27').Value = '  +9.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.631.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.48'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.996'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('E33').Value = '  -2.93%  '
$ws.Range('E34').Value = '  +2.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.23'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.56'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.387'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.58'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.86'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.78'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +14.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '148.74'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '150.67'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0550'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.616'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('E50').Value = '  +3.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0929'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.07%  '
